# Scheduled-runner update: refresh computed profit columns (H:N) across
# several leve-profit sheets with newly recalculated market-price figures.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 57147.5
$ws.Range("I68").Value = 15000
$ws.Range("J68").Value = 99295
$ws.Range("K68").Value = 15000
$ws.Range("L68").Value = 99295
$ws.Range("M68").Value = -14251
$ws.Range("N68").Value = -100793

$ws.Range("H71").Value = 57147.5
$ws.Range("I71").Value = 15000
$ws.Range("J71").Value = 99295
$ws.Range("K71").Value = 45000
$ws.Range("L71").Value = 297885
$ws.Range("M71").Value = -41256
$ws.Range("N71").Value = -305373

$ws.Range("H111").Value = 13749.375
$ws.Range("I111").Value = 25875
$ws.Range("J111").Value = 1623.75
$ws.Range("K111").Value = 77625
$ws.Range("L111").Value = 4871.25
$ws.Range("M111").Value = -74558
$ws.Range("N111").Value = -11005.25

$ws.Range("H138").Value = 21742506
$ws.Range("I138").Value = 4557.091
$ws.Range("J138").Value = 41668960
$ws.Range("K138").Value = 13671.273
$ws.Range("L138").Value = 125006880
$ws.Range("M138").Value = -8531.273000000001
$ws.Range("N138").Value = -125017160

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8726.846
$ws.Range("I2").Value = 10957.9
$ws.Range("J2").Value = 1290
$ws.Range("K2").Value = 10957.9
$ws.Range("L2").Value = 1290
$ws.Range("M2").Value = -10844.9
$ws.Range("N2").Value = -1516

$ws.Range("H45").Value = 1756.4445
$ws.Range("I45").Value = 1078.3636
$ws.Range("J45").Value = 4740
$ws.Range("K45").Value = 1078.3636
$ws.Range("L45").Value = 4740
$ws.Range("M45").Value = -701.3635999999999
$ws.Range("N45").Value = -5494

$ws.Range("H61").Value = 1856.4722
$ws.Range("I61").Value = 1067.9286
$ws.Range("J61").Value = 4616.375
$ws.Range("K61").Value = 1067.9286
$ws.Range("L61").Value = 4616.375
$ws.Range("M61").Value = -855.9286
$ws.Range("N61").Value = -5040.375

$ws.Range("H74").Value = 639.3514
$ws.Range("I74").Value = 489.73914
$ws.Range("J74").Value = 885.1429000000001
$ws.Range("K74").Value = 489.73914
$ws.Range("L74").Value = 885.1429000000001
$ws.Range("M74").Value = 384.26086
$ws.Range("N74").Value = -2633.1429

$ws.Range("H77").Value = 639.3514
$ws.Range("I77").Value = 489.73914
$ws.Range("J77").Value = 885.1429000000001
$ws.Range("K77").Value = 2448.6957
$ws.Range("L77").Value = 4425.7145
$ws.Range("M77").Value = 1919.3043
$ws.Range("N77").Value = -13161.7145

$ws.Range("H110").Value = 1036.1666
$ws.Range("I110").Value = 776.625
$ws.Range("J110").Value = 1555.25
$ws.Range("K110").Value = 776.625
$ws.Range("L110").Value = 1555.25
$ws.Range("M110").Value = 1268.375
$ws.Range("N110").Value = -5645.25

$ws.Range("H116").Value = 8726.846
$ws.Range("I116").Value = 10957.9
$ws.Range("J116").Value = 1290
$ws.Range("K116").Value = 10957.9
$ws.Range("L116").Value = 1290
$ws.Range("M116").Value = -8663.9
$ws.Range("N116").Value = -5878

$ws.Range("H132").Value = 2720.3872
$ws.Range("I132").Value = 2288.3684
$ws.Range("J132").Value = 3404.4167
$ws.Range("K132").Value = 6865.1052
$ws.Range("L132").Value = 10213.2501
$ws.Range("M132").Value = -4335.1052
$ws.Range("N132").Value = -15273.2501

$ws.Range("H136").Value = 1856.4722
$ws.Range("I136").Value = 1067.9286
$ws.Range("J136").Value = 4616.375
$ws.Range("K136").Value = 3203.7858
$ws.Range("L136").Value = 13849.125
$ws.Range("M136").Value = -653.7857999999997
$ws.Range("N136").Value = -18949.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8726.846
$ws.Range("I3").Value = 10957.9
$ws.Range("J3").Value = 1290
$ws.Range("K3").Value = 10957.9
$ws.Range("L3").Value = 1290
$ws.Range("M3").Value = -10843.9
$ws.Range("N3").Value = -1518

$ws.Range("H105").Value = 3015.1614
$ws.Range("I105").Value = 2778.8
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 2778.8
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -1031.8
$ws.Range("N105").Value = -7494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1145.4318
$ws.Range("I58").Value = 590.4516
$ws.Range("J58").Value = 2468.8462
$ws.Range("K58").Value = 590.4516
$ws.Range("L58").Value = 2468.8462
$ws.Range("M58").Value = -387.4516
$ws.Range("N58").Value = -2874.8462

$ws.Range("H105").Value = 655.1818
$ws.Range("I105").Value = 620.7
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 620.7
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 1126.3
$ws.Range("N105").Value = -4494

$ws.Range("H107").Value = 360
$ws.Range("I107").Value = 265
$ws.Range("J107").Value = 476.1111
$ws.Range("K107").Value = 265
$ws.Range("L107").Value = 476.1111
$ws.Range("M107").Value = 1655
$ws.Range("N107").Value = -4316.1111

$ws.Range("H132").Value = 3616.5715
$ws.Range("I132").Value = 2475.25
$ws.Range("J132").Value = 5138.3335
$ws.Range("K132").Value = 7425.75
$ws.Range("L132").Value = 15415.0005
$ws.Range("M132").Value = -4895.75
$ws.Range("N132").Value = -20475.0005

$ws.Range("H134").Value = 3875.1667
$ws.Range("I134").Value = 2613.5386
$ws.Range("J134").Value = 7155.4
$ws.Range("K134").Value = 7840.6158
$ws.Range("L134").Value = 21466.2
$ws.Range("M134").Value = -5305.6158
$ws.Range("N134").Value = -26536.2

$ws.Range("H136").Value = 1145.4318
$ws.Range("I136").Value = 590.4516
$ws.Range("J136").Value = 2468.8462
$ws.Range("K136").Value = 1771.3548
$ws.Range("L136").Value = 7406.5386
$ws.Range("M136").Value = 778.6451999999999
$ws.Range("N136").Value = -12506.5386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()

$ws.Range("H20").Value = 300
$ws.Range("I20").Value = 300
$ws.Range("J20").Value = 300
$ws.Range("K20").Value = 900
$ws.Range("L20").Value = 900
$ws.Range("M20").Value = -673
$ws.Range("N20").Value = -1354

$ws.Range("H103").Value = 145282.67
$ws.Range("I103").Value = 756.25
$ws.Range("J103").Value = 179288.88
$ws.Range("K103").Value = 2268.75
$ws.Range("L103").Value = 537866.64
$ws.Range("M103").Value = -1389.75
$ws.Range("N103").Value = -539624.64

$ws.Range("H113").Value = 371.65714
$ws.Range("I113").Value = 338.92307
$ws.Range("J113").Value = 391
$ws.Range("K113").Value = 1016.76921
$ws.Range("L113").Value = 1173
$ws.Range("M113").Value = 1153.23079
$ws.Range("N113").Value = -5513

$ws.Range("H131").Value = 2553.411
$ws.Range("I131").Value = 372.72726
$ws.Range("J131").Value = 2940.3064
$ws.Range("K131").Value = 1118.18178
$ws.Range("L131").Value = 8820.9192
$ws.Range("M131").Value = 3921.81822
$ws.Range("N131").Value = -18900.9192

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2651.5
$ws.Range("I102").Value = 2499.5
$ws.Range("J102").Value = 2803.5
$ws.Range("K102").Value = 2499.5
$ws.Range("L102").Value = 2803.5
$ws.Range("M102").Value = -877.5
$ws.Range("N102").Value = -6047.5

$ws.Range("H107").Value = 1113.3
$ws.Range("I107").Value = 1113.3
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1113.3
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 806.7
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3235.15
$ws.Range("I40").Value = 2650.5
$ws.Range("J40").Value = 3485.7144
$ws.Range("K40").Value = 2650.5
$ws.Range("L40").Value = 3485.7144
$ws.Range("M40").Value = -2514.5
$ws.Range("N40").Value = -3757.7144

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 13890162
$ws.Range("I132").Value = 16667493
$ws.Range("J132").Value = 3509
$ws.Range("K132").Value = 50002479
$ws.Range("L132").Value = 10527
$ws.Range("M132").Value = -49999949
$ws.Range("N132").Value = -15587
